$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link corrections for rows whose ranking order changed
$nameLinkUpdates = @(
    @{ Row = 35; Name = 'VeChain'; Link = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Row = 36; Name = 'FraxShare'; Link = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Row = 38; Name = 'Algorand'; Link = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Row = 39; Name = 'Hedera'; Link = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
)
foreach ($u in $nameLinkUpdates) {
    $ws.Range("B" + $u.Row).Value = $u.Name
    $ws.Range("C" + $u.Row).Value = $u.Link
}

# Updated price (column D, kept as text) and 1h volume change (column E) per row
$priceVolumeUpdates = @(
    @{ Row = 2; Price = '23.479.23'; Volume = '  +1.15%  ' }
    @{ Row = 3; Price = '1.646.72'; Volume = '  +2.84%  ' }
    @{ Row = 4; Price = '0.9992'; Volume = '  -0.27%  ' }
    @{ Row = 5; Price = '0.9989'; Volume = '  -0.26%  ' }
    @{ Row = 6; Price = '304.70'; Volume = '  +0.42%  ' }
    @{ Row = 7; Price = '0.3752'; Volume = '  -0.34%  ' }
    @{ Row = 8; Price = '51.98'; Volume = '  -0.27%  ' }
    @{ Row = 9; Price = '0.3632'; Volume = '  +0.07%  ' }
    @{ Row = 10; Price = '1.250'; Volume = '  -1.46%  ' }
    @{ Row = 11; Price = '0.08115'; Volume = '  -0.26%  ' }
    @{ Row = 12; Price = '0.9995'; Volume = '  -0.25%  ' }
    @{ Row = 13; Price = '22.89'; Volume = '  +0.75%  ' }
    @{ Row = 14; Price = '6.635'; Volume = '  +1.20%  ' }
    @{ Row = 15; Price = '0.00001269'; Volume = '  +1.90%  ' }
    @{ Row = 16; Price = '7.305'; Volume = '  -1.23%  ' }
    @{ Row = 17; Price = '1.635.12'; Volume = '  +2.22%  ' }
    @{ Row = 18; Price = '94.75'; Volume = '  +0.63%  ' }
    @{ Row = 19; Price = '0.06867'; Volume = '  -0.80%  ' }
    @{ Row = 20; Price = '18.21'; Volume = '  +0.77%  ' }
    @{ Row = 21; Price = '6.546'; Volume = '  +0.41%  ' }
    @{ Row = 22; Price = '0.9991'; Volume = '  -0.48%  ' }
    @{ Row = 23; Price = '23.490.99'; Volume = '  +1.23%  ' }
    @{ Row = 24; Price = '12.82'; Volume = '  -0.47%  ' }
    @{ Row = 25; Price = '3.149'; Volume = '  +3.87%  ' }
    @{ Row = 26; Price = '2.408'; Volume = '  -1.73%  ' }
    @{ Row = 27; Price = '21.21'; Volume = '  +0.26%  ' }
    @{ Row = 28; Price = '150.74'; Volume = '  +1.03%  ' }
    @{ Row = 29; Price = '5.299'; Volume = '  +0.81%  ' }
    @{ Row = 30; Price = '135.76'; Volume = '  +0.16%  ' }
    @{ Row = 31; Price = '2.285'; Volume = '  -3.86%  ' }
    @{ Row = 32; Price = '1.816.88'; Volume = '  +2.27%  ' }
    @{ Row = 33; Price = '6.828'; Volume = '  +1.99%  ' }
    @{ Row = 34; Price = '0.9558'; Volume = '  -0.47%  ' }
    @{ Row = 35; Price = '0.02812'; Volume = '  +2.78%  ' }
    @{ Row = 36; Price = '10.54'; Volume = '  +2.30%  ' }
    @{ Row = 37; Price = '6.252'; Volume = '  +2.79%  ' }
    @{ Row = 38; Price = '0.2535'; Volume = '  +0.49%  ' }
    @{ Row = 39; Price = '0.07282'; Volume = '  -2.45%  ' }
    @{ Row = 40; Price = '0.08835'; Volume = '  +0.56%  ' }
    @{ Row = 41; Price = '1.369'; Volume = '  -0.77%  ' }
    @{ Row = 42; Price = '0.7076'; Volume = '  +0.03%  ' }
    @{ Row = 43; Price = '12.49'; Volume = '  +0.85%  ' }
    @{ Row = 44; Price = '16.20'; Volume = '  +4.63%  ' }
    @{ Row = 45; Price = '0.6546'; Volume = '  +0.45%  ' }
    @{ Row = 46; Price = '2.337'; Volume = '  +1.21%  ' }
    @{ Row = 47; Price = '0.9984'; Volume = '  -0.21%  ' }
    @{ Row = 48; Price = '4.010'; Volume = '  +0.05%  ' }
    @{ Row = 49; Price = '0.07998'; Volume = '  +1.03%  ' }
    @{ Row = 50; Price = '129.20'; Volume = '  -2.12%  ' }
    @{ Row = 51; Price = '1.208'; Volume = '  +0.73%  ' }
)
foreach ($u in $priceVolumeUpdates) {
    $priceCell = $ws.Range("D" + $u.Row)
    # Force text storage so values such as "304.70" or "23.479.23" are not
    # reinterpreted as numbers/dates, matching the original inline-string cells.
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $u.Price
    $priceCell.Style = "Normal"
    $ws.Range("E" + $u.Row).Value = $u.Volume
}
